$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: shrink the explanation row height (112 -> 90) ---
$ws.Rows.Item(12).RowHeight = 90

# --- Row 15: new entry "Convert Sored Array to Binary Search Tree" ---
# Copy formats (cell styles) from row 14 so the new row reuses the same
# style indices rather than creating brand-new ones.
$ws.Range("A14:C14").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 142

$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Convert Sored Array to Binary Search Tree"

$explanation = "Recursion Technique.`n`nApproach:`n -- First perform binary search (ie find the mid)`n -- Make this mid as root node.`n -- Repeat the process on left sub tree and right sub tree."
$ws.Range("C15").Value = $explanation

# Base run font (matches the rest of the sheet's explanation column).
$full = $ws.Range("C15").Characters()
$full.Font.Name = "Times New Roman"
$full.Font.Size = 16

# First run: bold red heading "Recursion Technique.\n"
$headingLen = 21
$heading = $ws.Range("C15").Characters(1, $headingLen)
$heading.Font.Bold = $true
$heading.Font.Color = 255

# Second run: remaining normal (black) text
$restStart = $headingLen + 1
$restLen = $explanation.Length - $headingLen
$rest = $ws.Range("C15").Characters($restStart, $restLen)
$rest.Font.Color = 0

# --- Update dimension / selection bookkeeping ---
$ws.Range("A16").Select()
